$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G holds "K" (strikeouts), replacing old "Strike#" values with new K values
$ws.Range("G2").Value = 3
$ws.Range("G3").Value = 6
$ws.Range("G4").Value = 7
$ws.Range("G5").Value = 2
